# Added periodic & upfront related scenarios
# The "repaymentstrategy" input (row 17 on the ProductLoanInput sheet) is
# changed from "Mifos style" to "Penalties, Fees, Interest, Principal order",
# matching the new repayment-strategy scenario option that was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B17")
$cell.Select()
$cell.Value = "Penalties, Fees, Interest, Principal order"

# Left/top aligned (no wrap) to match the new repayment-strategy scenario text style
$cell.HorizontalAlignment = -4131  # xlLeft
$cell.VerticalAlignment = -4160    # xlTop
